$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 379.2
$ws.Range("I12").Value = 461
$ws.Range("K12").Value = 461
$ws.Range("M12").Value = -291

$ws.Range("H31").Value = 530.5
$ws.Range("I31").Value = 530.5
$ws.Range("K31").Value = 1591.5
$ws.Range("M31").Value = -1361.5

$ws.Range("H53").Value = 1237.5454
$ws.Range("I53").Value = 1580.9
$ws.Range("K53").Value = 1580.9
$ws.Range("M53").Value = -943.9000000000001

$ws.Range("H80").Value = 3640.9473
$ws.Range("I80").Value = 2500.1667
$ws.Range("J80").Value = 4167.4614
$ws.Range("K80").Value = 7500.500100000001
$ws.Range("L80").Value = 12502.3842
$ws.Range("M80").Value = -6502.500100000001
$ws.Range("N80").Value = -14498.3842

$ws.Range("H83").Value = 3640.9473
$ws.Range("I83").Value = 2500.1667
$ws.Range("J83").Value = 4167.4614
$ws.Range("K83").Value = 22501.5003
$ws.Range("L83").Value = 37507.1526
$ws.Range("M83").Value = -17509.5003
$ws.Range("N83").Value = -47491.1526

$ws.Range("H86").Value = 5512.9414
$ws.Range("J86").Value = 7401.3
$ws.Range("L86").Value = 7401.3
$ws.Range("N86").Value = -9647.3

$ws.Range("H89").Value = 5512.9414
$ws.Range("J89").Value = 7401.3
$ws.Range("L89").Value = 37006.5
$ws.Range("N89").Value = -48238.5

$ws.Range("H107").Value = 1711.875
$ws.Range("J107").Value = 1549.5
$ws.Range("L107").Value = 1549.5
$ws.Range("N107").Value = -5389.5

$ws.Range("H112").Value = 3442.8572
$ws.Range("J112").Value = 3500
$ws.Range("L112").Value = 10500
$ws.Range("N112").Value = -12716

$ws.Range("H132").Value = 34992.582
$ws.Range("I132").Value = 39689.594
$ws.Range("K132").Value = 119068.782
$ws.Range("M132").Value = -116538.782

$ws.Range("H138").Value = 2721.4146
$ws.Range("I138").Value = 2280.2
$ws.Range("J138").Value = 2975.9614
$ws.Range("K138").Value = 6840.599999999999
$ws.Range("L138").Value = 8927.8842
$ws.Range("M138").Value = -1700.599999999999
$ws.Range("N138").Value = -19207.8842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1610.4706
$ws.Range("I45").Value = 1493
$ws.Range("K45").Value = 1493
$ws.Range("M45").Value = -1116

$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 191.2
$ws.Range("I80").Value = 110
$ws.Range("J80").Value = 226
$ws.Range("K80").Value = 110
$ws.Range("L80").Value = 226
$ws.Range("M80").Value = 888
$ws.Range("N80").Value = -2222

$ws.Range("H83").Value = 191.2
$ws.Range("I83").Value = 110
$ws.Range("J83").Value = 226
$ws.Range("K83").Value = 550
$ws.Range("L83").Value = 1130
$ws.Range("M83").Value = 4442
$ws.Range("N83").Value = -11114

$ws.Range("H86").Value = 4029.1538
$ws.Range("I86").Value = 4701.75
$ws.Range("J86").Value = 3730.2222
$ws.Range("K86").Value = 4701.75
$ws.Range("L86").Value = 3730.2222
$ws.Range("M86").Value = -3578.75
$ws.Range("N86").Value = -5976.2222

$ws.Range("H89").Value = 4029.1538
$ws.Range("I89").Value = 4701.75
$ws.Range("J89").Value = 3730.2222
$ws.Range("K89").Value = 23508.75
$ws.Range("L89").Value = 18651.111
$ws.Range("M89").Value = -17892.75
$ws.Range("N89").Value = -29883.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 38249.43
$ws.Range("I122").Value = 2974.2
$ws.Range("K122").Value = 8922.599999999999
$ws.Range("M122").Value = -6472.599999999999

$ws.Range("H134").Value = 1526.2632
$ws.Range("I134").Value = 1646.8125
$ws.Range("K134").Value = 4940.4375
$ws.Range("M134").Value = -2405.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 14133.167
$ws.Range("J42").Value = 14133.167
$ws.Range("L42").Value = 42399.501
$ws.Range("N42").Value = -43467.501

$ws.Range("H100").Value = 5028
$ws.Range("J100").Value = 5028
$ws.Range("L100").Value = 15084
$ws.Range("N100").Value = -16706

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 26666.666
$ws.Range("J105").Value = 26666.666
$ws.Range("L105").Value = 79999.998
$ws.Range("N105").Value = -85241.998

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H139").Value = 2703.7
$ws.Range("I139").Value = 2282.125
$ws.Range("K139").Value = 6846.375
$ws.Range("M139").Value = -1706.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 154.94444
$ws.Range("I2").Value = 159.13333
$ws.Range("J2").Value = 134
$ws.Range("K2").Value = 159.13333
$ws.Range("L2").Value = 134
$ws.Range("M2").Value = -46.13333
$ws.Range("N2").Value = -360

$ws.Range("H80").Value = 2709.7
$ws.Range("I80").Value = 2869.6667
$ws.Range("J80").Value = 2469.75
$ws.Range("K80").Value = 2869.6667
$ws.Range("L80").Value = 2469.75
$ws.Range("M80").Value = -1871.6667
$ws.Range("N80").Value = -4465.75

$ws.Range("H83").Value = 2709.7
$ws.Range("I83").Value = 2869.6667
$ws.Range("J83").Value = 2469.75
$ws.Range("K83").Value = 14348.3335
$ws.Range("L83").Value = 12348.75
$ws.Range("M83").Value = -9356.3335
$ws.Range("N83").Value = -22332.75

$ws.Range("H93").Value = 15000
$ws.Range("I93").Value = 15000
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = -13128

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3486.7727
$ws.Range("I40").Value = 3160.9
$ws.Range("K40").Value = 3160.9
$ws.Range("M40").Value = -3024.9

$ws.Range("H87").Value = 88888
$ws.Range("J87").Value = 88888
$ws.Range("L87").Value = 88888
$ws.Range("N87").Value = -91134

$ws.Range("H90").Value = 88888
$ws.Range("J90").Value = 88888
$ws.Range("L90").Value = 266664
$ws.Range("N90").Value = -277896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10193.786
$ws.Range("J81").Value = 4599.5
$ws.Range("L81").Value = 9199
$ws.Range("N81").Value = -11321

$ws.Range("H84").Value = 10193.786
$ws.Range("J84").Value = 4599.5
$ws.Range("L84").Value = 45995
$ws.Range("N84").Value = -56603

$ws.Range("H126").Value = 3766
$ws.Range("I126").Value = 2100
$ws.Range("K126").Value = 6300
$ws.Range("M126").Value = -3830

Write-Host "Applied 35 row updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
